$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LSEG.L News Data")

$data = @(
    @("London Stock Exchange Group continues to underperform Intermediate Capital", "urn:newsml:newsroom:20231107:nNRAqll2kb:0", "NS:NEWMAR"),
    @("London Stock Exchange Group continues to underperform M&G", "urn:newsml:newsroom:20231107:nNRAqlk82m:0", "NS:NEWMAR"),
    @("London Stock Exchange Group's Earnings Growth Outperforms Hargreaves Lansdown", "urn:newsml:newsroom:20231107:nNRAqlk7t4:0", "NS:NEWMAR"),
    @("London Stock Exchange Group (LSEG: GBX8,314.0) - PRE-MARKET CONFIRMATORY BULLISH SIGNAL FROM ITS ADR MARKET Tuesday November 07, 2023 06:11 GMT", "urn:newsml:newsroom:20231107:nNRAqljhxn:0", "NS:GLOBUL"),
    @("UNITED KINGDOM PRE-MARKET BULLISH STOCK SNIPPETS LONDON STOCK EXCHANGE GROUP (LSEG: GBX8,314.0) Tuesday November 07, 2023 06:11 GMT", "urn:newsml:newsroom:20231107:nNRAqljii0:0", "NS:GLOBUL"),
    @("London Stock Exchange Group's Earnings Growth Lags 3I Group", "urn:newsml:newsroom:20231107:nNRAqljcyu:0", "NS:NEWMAR"),
    @("London Stock Exchange Group continues to outperform Schroders", "urn:newsml:newsroom:20231107:nNRAqljd3i:0", "NS:NEWMAR"),
    @("How London Stock Exchange Group's dividend compares to Intermediate Capital", "urn:newsml:newsroom:20231107:nNRAqljchb:0", "NS:NEWMAR"),
    @("Hargreaves Lansdown's Earnings Growth Lags London Stock Exchange Group", "urn:newsml:newsroom:20231107:nNRAqlif88:0", "NS:NEWMAR"),
    @("Close Brothers's Earnings Growth Lags London Stock Exchange Group", "urn:newsml:newsroom:20231107:nNRAqliewr:0", "NS:NEWMAR")
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
